# Updated TestData for Portugal Market
# - Adds a new "Portugal" sheet (cloned from the "Swiss" sheet layout)
# - Populates the Portugal-specific market name + ticket reference
# - Fixes up the Germany sheet's stale "select-all" selection

$wb = $excel.ActiveWorkbook

# --- Clone the Swiss sheet as the template for the new Portugal tab -------
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $swiss)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# --- Market-specific content ------------------------------------------------
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2459/T-2460/T2461"

# --- Column widths tuned for the new (shorter) strings ----------------------
$portugal.Columns.Item(1).ColumnWidth = 22.166666666666668
$portugal.Columns.Item(2).ColumnWidth = 26.166666666666668
$portugal.Columns.Item(3).ColumnWidth = 12.833333333333334
$portugal.Columns.Item(4).ColumnWidth = 11.666666666666666

# --- Row heights grow because column D now wraps onto two lines ------------
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# Row 16 carried over a stale custom height from the Swiss sheet; rebuild it
# (copying formats from row 15, which keeps its own height untouched) so it
# reverts back to the sheet's default row height.
$portugal.Range("A15").Copy()
$portugal.Rows.Item(16).Delete()
$portugal.Range("A16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$portugal.Range("A16").Value = "Printers"

# --- Active cell / tab state -------------------------------------------------
$portugal.Range("B4").Select()

# --- Germany sheet: selection had drifted to a full-sheet select -----------
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A1:D16").Select()

$portugal.Select()
